$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '88.691.55'
$ws.Range('E2').Value = '  +9.00%  '
$ws.Range('D3').Value = '3.341.22'
$ws.Range('E3').Value = '  +5.20%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '219.78'
$ws.Range('E5').Value = '  +5.68%  '
$ws.Range('D6').Value = '649.75'
$ws.Range('E6').Value = '  +2.75%  '
$ws.Range('D7').Value = '0.392'
$ws.Range('E7').Value = '  +33.16%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.608'
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('D10').Value = '3.343.62'
$ws.Range('E10').Value = '  +5.37%  '
$ws.Range('D11').Value = '0.589'
$ws.Range('E11').Value = '  +0.01%  '
$ws.Range('D12').Value = '0.0000283'
$ws.Range('E12').Value = '  +7.91%  '
$ws.Range('E13').Value = '  +1.74%  '
$ws.Range('D14').Value = '35.36'
$ws.Range('E14').Value = '  +10.51%  '
$ws.Range('D15').Value = '3.952.96'
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').Value = '88.520.65'
$ws.Range('E17').Value = '  +8.73%  '
$ws.Range('D18').Value = '3.346.53'
$ws.Range('E18').Value = '  +5.55%  '
$ws.Range('D19').Value = '14.67'
$ws.Range('E19').Value = '  +3.28%  '
$ws.Range('D20').Value = '3.17'
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('D21').Value = '9.74'
$ws.Range('E21').Value = '  +5.74%  '
$ws.Range('D22').Value = '458.44'
$ws.Range('E22').Value = '  +4.40%  '
$ws.Range('E23').Value = '  +6.28%  '
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').Value = '7.43'
$ws.Range('E24').Value = '  +4.10%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = '5.61'
$ws.Range('E25').Value = '  +8.12%  '
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').Value = '12.86'
$ws.Range('E26').Value = '  +14.58%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '3.514.34'
$ws.Range('E27').Value = '  +5.11%  '
$ws.Range('B28').Value = 'Litecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D28').Value = '78.82'
$ws.Range('E28').Value = '  +2.47%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0000130'
$ws.Range('E29').Value = '  +3.24%  '
$ws.Range('B30').Value = 'Cronos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D30').Value = '0.199'
$ws.Range('E30').Value = '  +42.81%  '
$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').Value = '9.40'
$ws.Range('E32').Value = '  +2.85%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '600.80'
$ws.Range('E33').Value = '  +7.09%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '1.61'
$ws.Range('E34').Value = '  +6.25%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').Value = '0.991'
$ws.Range('E35').Value = '  -0.91%  '
$ws.Range('B36').Value = 'PancakeSwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D36').Value = '2.15'
$ws.Range('E36').Value = '  +5.61%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '7.31'
$ws.Range('E37').Value = '  +22.61%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '0.147'
$ws.Range('E38').Value = '  -3.57%  '
$ws.Range('B39').Value = 'EthereumClassic'
$ws.Range('C39').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D39').Value = '23.42'
$ws.Range('E39').Value = '  +1.15%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '2.19'
$ws.Range('E40').Value = '  +7.85%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').Value = '0.421'
$ws.Range('E41').Value = '  +1.65%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = '21.88'
$ws.Range('E42').Value = '  +5.43%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '3.11'
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').Value = '159.10'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '190.59'
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('B48').Value = 'ImmutableX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D48').Value = '1.44'
$ws.Range('E48').Value = '  +6.99%  '
$ws.Range('B49').Value = 'OKB'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D49').Value = '46.55'
$ws.Range('E49').Value = '  +4.50%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.787'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').Value = '4.43'
$ws.Range('E51').Value = '  +3.77%  '
